$d = $word.ActiveDocument

# 1. Replace the hyperlink display text (by-nc -> by-sa)
$r = $d.Content.Find.Execute("https://creativecommons.org/licenses/by-nc/4.0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://creativecommons.org/licenses/by-sa/4.0", 2)

# 2. Replace the inline licence text "CC BY-NC" -> "CC BY-SA"
$r2 = $d.Content.Find.Execute("CC BY-NC", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CC BY-SA", 2)

# 3. Update the actual hyperlink relationship target
foreach ($hl in $d.Hyperlinks) {
    if ($hl.Address -eq "https://creativecommons.org/licenses/by-nc/4.0") {
        $hl.Address = "https://creativecommons.org/licenses/by-sa/4.0"
    }
}
